$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.878.76"
$ws.Range("E2").Value = "  +1.24%  "

$ws.Range("D3").Value = "1.840.80"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.006"
$ws.Range("E4").Value = "  +0.26%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "309.35"
$ws.Range("E5").Value = "  +1.09%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.006"
$ws.Range("E6").Value = "  +0.33%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4718"
$ws.Range("E7").Value = "  +3.72%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3648"
$ws.Range("E8").Value = "  +1.50%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07151"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9192"
$ws.Range("E10").Value = "  +3.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.50"

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07606"
$ws.Range("E12").Value = "  -1.37%  "

$ws.Range("D13").Value = "1.787.65"
$ws.Range("E13").Value = "  -1.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.278"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.389"
$ws.Range("E15").Value = "  +1.48%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "87.77"
$ws.Range("E16").Value = "  +1.35%  "

$ws.Range("E17").Value = "  +0.43%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008624"
$ws.Range("E18").Value = "  +0.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.006"
$ws.Range("E19").Value = "  +0.31%  "

$ws.Range("D20").Value = "26.913.57"
$ws.Range("E20").Value = "  +1.25%  "

$ws.Range("E21").Value = "  +2.38%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.008"
$ws.Range("E22").Value = "  +0.93%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.60"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.925"
$ws.Range("E24").Value = "  -0.20%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "151.57"
$ws.Range("E25").Value = "  -0.07%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "18.17"
$ws.Range("E26").Value = "  +2.10%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.006"
$ws.Range("E27").Value = "  -0.51%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.04"
$ws.Range("E28").Value = "  +1.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.08818"
$ws.Range("E30").Value = "  +1.11%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.232"
$ws.Range("E31").Value = "  +3.43%  "

$ws.Range("E32").Value = "  +5.25%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7409"
$ws.Range("E33").Value = "  +0.33%  "

$ws.Range("B34").Value = "RenderToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.748"
$ws.Range("E34").Value = "  +1.06%  "

$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.472"
$ws.Range("E35").Value = "  +0.87%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.087"
$ws.Range("E36").Value = "  +1.54%  "

$ws.Range("E37").Value = "  +0.39%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05237"
$ws.Range("E38").Value = "  +3.28%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.963"
$ws.Range("E39").Value = "  +1.77%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5179"
$ws.Range("E40").Value = "  +1.95%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.935"
$ws.Range("E41").Value = "  +2.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1510"
$ws.Range("E42").Value = "  +0.28%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.144"
$ws.Range("E43").Value = "  +1.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.41"
$ws.Range("E44").Value = "  +4.77%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4698"
$ws.Range("E45").Value = "  +0.30%  "

$ws.Range("E46").Value = "  +0.41%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.81"
$ws.Range("E47").Value = "  +2.54%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.591"
$ws.Range("E48").Value = "  +1.56%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.27"
$ws.Range("E49").Value = "  +2.51%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06030"
$ws.Range("E50").Value = "  +0.44%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8849"
$ws.Range("E51").Value = "  +4.46%  "
